$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 474.51614
$ws.Range("J17").Value = 474.51614
$ws.Range("L17").Value = 1423.54842
$ws.Range("N17").Value = -1759.54842

# Row 38
$ws.Range("H38").Value = 370.69232
$ws.Range("I38").Value = 76.583336
$ws.Range("K38").Value = 229.750008
$ws.Range("M38").Value = 142.249992

# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 111
$ws.Range("H111").Value = 3832.6667
$ws.Range("I111").Value = 1938
$ws.Range("K111").Value = 5814
$ws.Range("M111").Value = -2747

# Row 113
$ws.Range("H113").Value = 2886.111
$ws.Range("I113").Value = 2496.6667
$ws.Range("J113").Value = 3080.8333
$ws.Range("K113").Value = 2496.6667
$ws.Range("L113").Value = 3080.8333
$ws.Range("M113").Value = 757.3332999999998
$ws.Range("N113").Value = -9588.8333

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12392.2
$ws.Range("I32").Value = 14523.62
$ws.Range("J32").Value = 5287.467
$ws.Range("K32").Value = 14523.62
$ws.Range("L32").Value = 5287.467
$ws.Range("M32").Value = -14236.62
$ws.Range("N32").Value = -5861.467

# Row 63
$ws.Range("H63").Value = 4200
$ws.Range("I63").Value = 3840
$ws.Range("J63").Value = 6000
$ws.Range("K63").Value = 3840
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -3154
$ws.Range("N63").Value = -7372

# Row 66
$ws.Range("H66").Value = 4200
$ws.Range("I66").Value = 3840
$ws.Range("J66").Value = 6000
$ws.Range("K66").Value = 19200
$ws.Range("L66").Value = 30000
$ws.Range("M66").Value = -15768
$ws.Range("N66").Value = -36864

# Row 80
$ws.Range("H80").Value = 35633.332
$ws.Range("J80").Value = 35633.332
$ws.Range("L80").Value = 35633.332
$ws.Range("N80").Value = -37629.332

# Row 83
$ws.Range("H83").Value = 35633.332
$ws.Range("J83").Value = 35633.332
$ws.Range("L83").Value = 106899.996
$ws.Range("N83").Value = -116883.996

# Row 97
$ws.Range("H97").Value = 840.4286
$ws.Range("I97").Value = 810.6
$ws.Range("K97").Value = 810.6
$ws.Range("M97").Value = -314.6

$ws = $wb.Worksheets.Item("BSM")
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# Row 82
$ws.Range("H82").Value = 42281.25
$ws.Range("I82").Value = 42281.25
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 42281.25
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -41898.25
$ws.Range("N82").ClearContents()

# Row 85
$ws.Range("H85").Value = 42281.25
$ws.Range("I85").Value = 42281.25
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 42281.25
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -40955.25
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1446.0869
$ws.Range("I16").Value = 1331.2
$ws.Range("J16").Value = 1534.4615
$ws.Range("K16").Value = 1331.2
$ws.Range("L16").Value = 1534.4615
$ws.Range("M16").Value = -1044.2
$ws.Range("N16").Value = -2108.4615

# Row 31
$ws.Range("H31").Value = 3607.4092
$ws.Range("I31").Value = 1786.6471
$ws.Range("J31").Value = 9798
$ws.Range("K31").Value = 1786.6471
$ws.Range("L31").Value = 9798
$ws.Range("M31").Value = -1491.6471
$ws.Range("N31").Value = -10388

# Row 34
$ws.Range("H34").Value = 3607.4092
$ws.Range("I34").Value = 1786.6471
$ws.Range("J34").Value = 9798
$ws.Range("K34").Value = 1786.6471
$ws.Range("L34").Value = 9798
$ws.Range("M34").Value = -1584.6471
$ws.Range("N34").Value = -10202

# Row 113
$ws.Range("H113").Value = 1446.0869
$ws.Range("I113").Value = 1331.2
$ws.Range("J113").Value = 1534.4615
$ws.Range("K113").Value = 1331.2
$ws.Range("L113").Value = 1534.4615
$ws.Range("M113").Value = 838.8
$ws.Range("N113").Value = -5874.461499999999

# Row 132
$ws.Range("H132").Value = 1963.0256
$ws.Range("I132").Value = 1491.1111
$ws.Range("J132").Value = 3024.8333
$ws.Range("K132").Value = 4473.3333
$ws.Range("L132").Value = 9074.499899999999
$ws.Range("M132").Value = -1943.3333
$ws.Range("N132").Value = -14134.4999

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 4760.8
$ws.Range("I5").Value = 5701
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 17103
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -16991
$ws.Range("N5").Value = -3224

# Row 33
$ws.Range("H33").Value = 1317.5834
$ws.Range("J33").Value = 1578.875
$ws.Range("L33").Value = 9473.25
$ws.Range("N33").Value = -10039.25

# Row 122
$ws.Range("H122").Value = 907.4231
$ws.Range("I122").Value = 802
$ws.Range("K122").Value = 7218
$ws.Range("M122").Value = -4768

# Row 132
$ws.Range("H132").Value = 1782.0741
$ws.Range("I132").Value = 1204.6
$ws.Range("J132").Value = 2121.7646
$ws.Range("K132").Value = 10841.4
$ws.Range("L132").Value = 19095.8814
$ws.Range("M132").Value = -8311.4
$ws.Range("N132").Value = -24155.8814

# Row 135
$ws.Range("H135").Value = 4760.8
$ws.Range("I135").Value = 5701
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 51309
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -48774
$ws.Range("N135").Value = -14070

# Row 139
$ws.Range("H139").Value = 1857.0333
$ws.Range("I139").Value = 1276.0526
$ws.Range("J139").Value = 2860.5454
$ws.Range("K139").Value = 3828.1578
$ws.Range("L139").Value = 8581.636200000001
$ws.Range("M139").Value = 1311.8422
$ws.Range("N139").Value = -18861.6362

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3336.3635
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002

# Row 83
$ws.Range("H83").Value = 3336.3635
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("M83").Value = -10008

# Row 123
$ws.Range("H123").Value = 9000.25
$ws.Range("J123").Value = 9000.25
$ws.Range("L123").Value = 9000.25
$ws.Range("N123").Value = -13900.25

$ws = $wb.Worksheets.Item("LTW")
# Row 75
$ws.Range("H75").Value = 26000
$ws.Range("J75").Value = 26000
$ws.Range("L75").Value = 26000
$ws.Range("N75").Value = -27872

# Row 78
$ws.Range("H78").Value = 26000
$ws.Range("J78").Value = 26000
$ws.Range("L78").Value = 78000
$ws.Range("N78").Value = -87360

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 66265.8
$ws.Range("J46").Value = 66265.8
$ws.Range("L46").Value = 66265.8
$ws.Range("N46").Value = -66727.8

# Row 81
$ws.Range("H81").Value = 89206.30499999999
$ws.Range("I81").Value = 161826
$ws.Range("J81").Value = 4483.3335
$ws.Range("K81").Value = 323652
$ws.Range("L81").Value = 8966.666999999999
$ws.Range("M81").Value = -322591
$ws.Range("N81").Value = -11088.667

# Row 84
$ws.Range("H84").Value = 89206.30499999999
$ws.Range("I84").Value = 161826
$ws.Range("J84").Value = 4483.3335
$ws.Range("K84").Value = 1618260
$ws.Range("L84").Value = 44833.335
$ws.Range("M84").Value = -1612956
$ws.Range("N84").Value = -55441.335

# Row 125
$ws.Range("H125").Value = 60480
$ws.Range("J125").Value = 60480
$ws.Range("L125").Value = 60480
$ws.Range("N125").Value = -70320

# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 132
$ws.Range("H132").Value = 4364.0713
$ws.Range("I132").Value = 3586.1428
$ws.Range("J132").Value = 5142
$ws.Range("K132").Value = 10758.4284
$ws.Range("L132").Value = 15426
$ws.Range("M132").Value = -8228.428400000001
$ws.Range("N132").Value = -20486

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Row 134
$ws.Range("H134").Value = 66265.8
$ws.Range("J134").Value = 66265.8
$ws.Range("L134").Value = 198797.4
$ws.Range("N134").Value = -203867.4
